# Sablefish_Inputs.xlsx - "updating vonB sd paramterization"
#
# The single von Bertalanffy SD parameter ("vonB_sd", value 7/7) is split
# into two separate parameters:
#   - vonB_sd_1 = 3 / 3   (new row, takes over row 7)
#   - vonB_sd_2 = 7 / 7   (keeps the old vonB_sd value, shifts to row 8)
# and "wl_sd" (previously row 8, value 2/2) moves down to row 9 with a new
# value of 0.2 / 0.2.
#
# Various sheet-view/selection bookkeeping (active sheet/cell) also shifts
# around as a side effect of the editor session that made the change.

$wb = $excel.ActiveWorkbook

# --- Bins: the last bin edge (row 32, B32=100) was removed ---
$wsBins = $wb.Worksheets.Item("Bins")
$wsBins.Rows.Item(32).Delete() | Out-Null
$wsBins.Activate() | Out-Null
$wsBins.Range("E31").Select() | Out-Null

# --- Growth_Param: re-parameterize vonB_sd into vonB_sd_1 / vonB_sd_2,
#     and give wl_sd its own row with a new value ---
$wsGrowth = $wb.Worksheets.Item("Growth_Param")

# Write row 8 (vonB_sd_2) first so it lands earlier in the shared-string
# table than vonB_sd_1, matching the source ordering.
$wsGrowth.Range("A8").Value = 7
$wsGrowth.Range("B8").Value = 7
$wsGrowth.Range("C8").Value = "vonB_sd_2"

$wsGrowth.Range("A7").Value = 3
$wsGrowth.Range("B7").Value = 3
$wsGrowth.Range("C7").Value = "vonB_sd_1"

$wsGrowth.Range("A9").Value = 0.2
$wsGrowth.Range("B9").Value = 0.2
$wsGrowth.Range("C9").Value = "wl_sd"

$wsGrowth.Activate() | Out-Null
$wsGrowth.Range("D15").Select() | Out-Null

# --- Controls: becomes the active sheet/tab when the file was saved ---
$wsControls = $wb.Worksheets.Item("Controls")
$wsControls.Activate() | Out-Null
$wsControls.Range("B3").Select() | Out-Null
